$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# Add a new slide after slide 1, using the same layout as slide 1 (Title Slide)
$s2 = $p.Slides.Add(2, 1)

# Remove the default placeholder shapes (Title, Subtitle) that came with the layout
while ($s2.Shapes.Count -gt 0) {
    $s2.Shapes.Item(1).Delete()
}

# Copy the legend table (first table) from slide 1 and paste it onto the new slide
$srcTable = $s1.Shapes.Item(2)
$srcTable.Copy()
$s2.Shapes.Paste() | Out-Null
$newShape = $s2.Shapes.Item($s2.Shapes.Count)

$tbl = $newShape.Table

# The source table has 6 rows; keep only the first two (header legend + one data row)
$tbl.Rows.Item(6).Delete()
$tbl.Rows.Item(5).Delete()
$tbl.Rows.Item(4).Delete()
$tbl.Rows.Item(3).Delete()

# Update the second row's values for the finished breadth-first-search state
$tbl.Cell(2,1).Shape.TextFrame.TextRange.Text = "4"
$tbl.Cell(2,2).Shape.TextFrame.TextRange.Text = "2"
$tbl.Cell(2,3).Shape.TextFrame.TextRange.Text = "3"
$tbl.Cell(2,4).Shape.TextFrame.TextRange.Text = "1"
$tbl.Cell(2,5).Shape.TextFrame.TextRange.Text = "0"

# Reposition the table slightly (matches the nudged position on the new slide)
$newShape.Left = 856.1146240234375
$newShape.Top = 36.56248092651367
